# Proyecto Grupal / Datos.xlsx
# Adds column M ("promedio"/ANOVA helper column) = AVERAGE(K,E,B) for each
# data row (5..26), matching the Python-driven ANOVA analysis described in
# the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M: average of the three "tiempo" columns (K, E, B) ------
$dataRange = $ws.Range("M5:M26")
$dataRange.NumberFormat = "0.0"
$dataRange.Formula = "=AVERAGE(K5,E5,B5)"

# --- View state: mirror the author's re-selection of the new column -----
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$excel.Goto($ws.Range("M5:M26"))
